$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# 1) H5 gets the "Next" marker that used to live in G6 (same text + style).
$ws.Range("G6").Copy($ws.Range("H5")) | Out-Null

# 2) G6 becomes a date cell, picking up the plain date style already used
#    by F5/G5 (fill-only highlight, no special font) and the new date value.
$ws.Range("G5").Copy() | Out-Null
$ws.Range("G6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G6").Value = "2025-08-27"

# 3) F5 switches from the (now unused) bold/red "Next" date style to the
#    plain date style, keeping its existing value.
$ws.Range("G5").Copy() | Out-Null
$ws.Range("F5").PasteSpecial($xlPasteFormats) | Out-Null

# 4) O6 gets a date value (format/style already correct).
$ws.Range("O6").Value = "2025-08-27"

# 5) A8 is a new, empty, formatted cell matching A6's style.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A8").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# 6) Leave the selection where the author last clicked.
$ws.Range("J9").Select() | Out-Null
